$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the IT8 test (rows 9-12): move the "T" (top module) marker to the
# correct column for each step, and set the previous "T" cell (column B) to "X".
$ws.Range("B9").Value = "X"
$ws.Range("G9").Value = "T"

$ws.Range("B10").Value = "X"
$ws.Range("F10").Value = "T"

$ws.Range("B11").Value = "X"
$ws.Range("D11").Value = "T"

$ws.Range("B12").Value = "X"
$ws.Range("E12").Value = "T"

# Update the active cell selection to match the latest edit.
$ws.Range("E12").Select()
